# Component 3 deck - "PLENARY- COMPLETE THE MISSING GAPS" slide (slide 16):
# the cash-flow recap table had its Table Style switched (via the Tables
# Design ribbon) from the old built-in style to a new one.
#
#   old style id: {6201D1A8-7C20-4C7A-8E2F-EC44211F9D0D}
#   new style id: {802B5B32-C7BE-46BF-9F80-0C8C666AAB4F}

$p = $ppt.ActivePresentation

$targetStyleId = "{802B5B32-C7BE-46BF-9F80-0C8C666AAB4F}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shp = $slide.Shapes.Item($j)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            # Table styles are identified by a GUID and can't be assigned
            # through the Style property directly - PowerPoint applies them
            # via ApplyStyle (mirrors picking a style from the Table Styles
            # gallery on the ribbon).
            $tbl.ApplyStyle($targetStyleId)
        }
    }
}
